$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend header row (row 1): add P1=14, Q1=15 with same style as O1 (bold/border/centered)
$ws.Range("O1").Copy() | Out-Null
$ws.Range("P1:Q1").PasteSpecial(-4122) | Out-Null
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Update data rows 2-25: columns I-O get new values, and new columns P,Q are added.
# New values for columns I..Q (9 values), repeated per row:
$newValues = @(2, 2, 1, 2, 2, 2, 1, 2, 2)
$cols = @("I", "J", "K", "L", "M", "N", "O", "P", "Q")

for ($r = 2; $r -le 25; $r++) {
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $r).Value = $newValues[$i]
    }
}
